$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 48 of portfolio data for 2025-10-02.
# Column A holds a plain text date label (matching the style of the
# existing rows, which store dates as literal text rather than Excel
# date serials). Prefixing with an apostrophe forces text entry so the
# date-like string "2025-10-02" is not auto-converted to a date value.
$ws.Range("A48").Value = "'2025-10-02"
$ws.Range("B48").Value = 55.20000076293945
$ws.Range("C48").Value = 718
$ws.Range("D48").Value = 329.4500122070312
